$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.292504191398621
$ws.Range("B1").Value = 4.988254070281982
$ws.Range("C1").Value = 3.196216344833374
$ws.Range("D1").Value = 1.72171413898468
$ws.Range("E1").Value = 1.30008602142334
